# Trade #51 closed at 2026-02-17 13:28:49 - unknown UNKNOWN +0.000%
#
# This script:
#   1) Updates the "Summary" sheet roll-up metrics (B3:B9).
#   2) Updates the "Strategy Status" sheet row for MarketMaking (C4:G4).
#   3) Appends the new trade (#51) as row 52 on both the "All Trades" and
#      "MarketMaking" sheets (they mirror the same trade log).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) Summary sheet
# ---------------------------------------------------------------------
$summary = $wb.Worksheets.Item("Summary")
$summary.Range("B3").Value = 1197.42   # Current Capital
$summary.Range("B4").Value = -2.57    # Total P&L $
$summary.Range("B5").Value = -1.01    # Total P&L %
$summary.Range("B6").Value = 51       # Total Trades
$summary.Range("B8").Value = 27       # Losing Trades
$summary.Range("B9").Value = 37.25    # Win Rate %

# ---------------------------------------------------------------------
# 2) Strategy Status sheet - MarketMaking row (row 4)
# ---------------------------------------------------------------------
$status = $wb.Worksheets.Item("Strategy Status")
$status.Range("C4").Value = 97.42    # Capital
$status.Range("D4").Value = 51       # Trades
$status.Range("E4").Value = -2.57    # P&L $
$status.Range("F4").Value = -2.58    # P&L %
$status.Range("G4").Value = 37.25    # Win Rate %

# ---------------------------------------------------------------------
# 3) Append trade #51 as row 52 to both trade-log sheets
# ---------------------------------------------------------------------
function Add-TradeRow($ws) {
    $ws.Cells.Item(52, 1).Value = 51

    # A leading apostrophe forces text entry (mirrors typing '2026-02-17
    # into Excel) so the ISO-looking date string isn't auto-converted to
    # a date serial; resetting the style afterwards drops the implicit
    # "quote prefix" formatting Excel would otherwise remember.
    $ws.Cells.Item(52, 2).Value = "'2026-02-17"
    $ws.Cells.Item(52, 2).Style = "Normal"

    $ws.Cells.Item(52, 3).Value = "13:28:43"
    $ws.Cells.Item(52, 4).Value = "MarketMaking"
    $ws.Cells.Item(52, 5).Value = "UP"
    $ws.Cells.Item(52, 6).Value = 0.98
    $ws.Cells.Item(52, 7).Value = 0.9
    $ws.Cells.Item(52, 8).Value = "CLOSED"
    $ws.Cells.Item(52, 9).Value = -8.1633
    $ws.Cells.Item(52, 10).Value = -0.08
    $ws.Cells.Item(52, 11).Value = 97.42
    $ws.Cells.Item(52, 12).Value = 0
    $ws.Cells.Item(52, 13).Value = 0
    $ws.Cells.Item(52, 14).Value = 0.6
    $ws.Cells.Item(52, 15).Value = "Normal spread capture: 19600 bps"
    $ws.Cells.Item(52, 16).Value = "early_exit"
    $ws.Cells.Item(52, 17).Value = 0.13
}

Add-TradeRow($wb.Worksheets.Item("All Trades"))
Add-TradeRow($wb.Worksheets.Item("MarketMaking"))
